$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link (columns B/C) are always plain text, assign directly.
# Price/Volume (columns D/E) are numeric-looking strings in the source data
# (e.g. "90.872.73", "1.00", "8.64", "  +1.85%  ") that must stay literal
# text rather than being auto-parsed into numbers/percentages by Excel -
# a leading apostrophe forces text entry, exactly like typing it by hand.

# Row 2
$ws.Range('D2').Value = "'90.872.73"
$ws.Range('E2').Value = "'  +1.85%  "

# Row 3
$ws.Range('D3').Value = "'3.210.86"
$ws.Range('E3').Value = "'  -2.44%  "

# Row 4
$ws.Range('E4').Value = "'  +0.27%  "

# Row 5
$ws.Range('D5').Value = "'215.24"
$ws.Range('E5').Value = "'  +0.50%  "

# Row 6
$ws.Range('D6').Value = "'617.38"
$ws.Range('E6').Value = "'  -2.01%  "

# Row 7
$ws.Range('D7').Value = "'0.387"
$ws.Range('E7').Value = "'  +0.07%  "

# Row 8
$ws.Range('D8').Value = "'0.695"
$ws.Range('E8').Value = "'  +0.59%  "

# Row 9
$ws.Range('E9').Value = "'  +0.15%  "

# Row 10
$ws.Range('D10').Value = "'3.197.34"
$ws.Range('E10').Value = "'  -2.67%  "

# Row 11
$ws.Range('D11').Value = "'0.575"
$ws.Range('E11').Value = "'  -0.85%  "

# Row 12
$ws.Range('D12').Value = "'0.178"
$ws.Range('E12').Value = "'  -4.78%  "

# Row 13
$ws.Range('D13').Value = "'0.0000257"
$ws.Range('E13').Value = "'  -2.73%  "

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = "'3.862.76"
$ws.Range('E14').Value = "'  -0.63%  "

# Row 15
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = "'90.870.60"
$ws.Range('E15').Value = "'  +2.14%  "

# Row 16
$ws.Range('D16').Value = "'33.34"
$ws.Range('E16').Value = "'  -2.65%  "

# Row 17
$ws.Range('D17').Value = "'5.23"
$ws.Range('E17').Value = "'  -2.88%  "

# Row 18
$ws.Range('D18').Value = "'3.246.15"
$ws.Range('E18').Value = "'  -1.43%  "

# Row 19
$ws.Range('D19').Value = "'3.24"
$ws.Range('E19').Value = "'  +3.97%  "

# Row 20
$ws.Range('D20').Value = "'13.57"
$ws.Range('E20').Value = "'  -4.29%  "

# Row 21
$ws.Range('D21').Value = "'436.15"
$ws.Range('E21').Value = "'  -0.41%  "

# Row 22
$ws.Range('B22').Value = 'PEPE'
$ws.Range('C22').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D22').Value = "'0.0000187"
$ws.Range('E22').Value = "'  +39.42%  "

# Row 23
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'8.64"
$ws.Range('E23').Value = "'  -3.02%  "

# Row 24
$ws.Range('D24').Value = "'5.15"
$ws.Range('E24').Value = "'  -4.92%  "

# Row 25
$ws.Range('D25').Value = "'5.22"
$ws.Range('E25').Value = "'  -0.61%  "

# Row 26
$ws.Range('D26').Value = "'11.85"
$ws.Range('E26').Value = "'  -4.33%  "

# Row 27
$ws.Range('D27').Value = "'3.448.15"
$ws.Range('E27').Value = "'  -0.05%  "

# Row 28
$ws.Range('D28').Value = "'76.03"
$ws.Range('E28').Value = "'  -1.55%  "

# Row 29
$ws.Range('E29').Value = "'  -0.07%  "

# Row 30
$ws.Range('D30').Value = "'0.170"
$ws.Range('E30').Value = "'  -11.59%  "

# Row 31
$ws.Range('D31').Value = "'1.00"
$ws.Range('E31').Value = "'  +0.06%  "

# Row 32
$ws.Range('E32').Value = "'  +39.58%  "

# Row 33
$ws.Range('D33').Value = "'8.51"
$ws.Range('E33').Value = "'  -4.26%  "

# Row 34
$ws.Range('D34').Value = "'537.89"
$ws.Range('E34').Value = "'  -6.73%  "

# Row 35
$ws.Range('D35').Value = "'6.87"
$ws.Range('E35').Value = "'  -5.69%  "

# Row 36
$ws.Range('D36').Value = "'1.89"
$ws.Range('E36').Value = "'  -4.25%  "

# Row 37
$ws.Range('D37').Value = "'1.24"
$ws.Range('E37').Value = "'  -11.32%  "

# Row 38
$ws.Range('D38').Value = "'22.23"
$ws.Range('E38').Value = "'  -2.15%  "

# Row 39
$ws.Range('D39').Value = "'22.38"
$ws.Range('E39').Value = "'  +2.56%  "

# Row 40
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = "'  +0.36%  "

# Row 41
$ws.Range('D41').Value = "'0.126"
$ws.Range('E41').Value = "'  -9.31%  "

# Row 42
$ws.Range('D42').Value = "'0.381"
$ws.Range('E42').Value = "'  -4.88%  "

# Row 43
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = "'  -0.08%  "

# Row 44
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = "'1.93"
$ws.Range('E44').Value = "'  -5.16%  "

# Row 45
$ws.Range('D45').Value = "'148.73"
$ws.Range('E45').Value = "'  -3.66%  "

# Row 46
$ws.Range('D46').Value = "'44.88"
$ws.Range('E46').Value = "'  -0.39%  "

# Row 47
$ws.Range('D47').Value = "'174.69"
$ws.Range('E47').Value = "'  -3.61%  "

# Row 48
$ws.Range('D48').Value = "'0.124"
$ws.Range('E48').Value = "'  -1.52%  "

# Row 49
$ws.Range('D49').Value = "'1.24"
$ws.Range('E49').Value = "'  -5.19%  "

# Row 50
$ws.Range('D50').Value = "'0.621"
$ws.Range('E50').Value = "'  -1.12%  "

# Row 51
$ws.Range('D51').Value = "'4.11"
$ws.Range('E51').Value = "'  -3.40%  "

